# Generate Report for Handoff
# Refresh the localization-status report: the handback run picked up a new
# source guid (9903108f-ca62-4da0-b928-b03a1d35df4d) and fresh xliff/handoff
# hashes+timestamps. Update the Overview + per-locale sheets to match, and
# clear the now-stale "Latest Target/Handback File" columns (the de-de
# handback hasn't landed yet, so its handback datetime resets to the zero
# date and its target/handback file links go away).

$wb = $excel.ActiveWorkbook

$oldGuid = "3985b3b1-8d61-43bb-9f53-b88b3951a67d"
$newGuid = "9903108f-ca62-4da0-b928-b03a1d35df4d"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = "e2e\$newGuid.md"
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}
$ws.Range("G2").Value = "2016-08-30 13:04:37"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newGuid.md"
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid.md"
    }
}
$ws.Range("G2").Value = "$newGuid.c0a71245d14ecb70dbc22e5c0b85f61e32a10d17.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-30 13:04:31"

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$I$2') {
        $hl.Delete()
    }
}
$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"

$ws.Range("J2").Value = ""
$ws.Range("J2").Style = "Normal"

$ws.Range("K2").Value = "0001-01-01 00:00:00"

$ws.Columns.Item(9).ColumnWidth = 17.83
$ws.Columns.Item(10).ColumnWidth = 20.85

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newGuid.md"
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid.md"
    }
}
$ws.Range("G2").Value = "$newGuid.c0a71245d14ecb70dbc22e5c0b85f61e32a10d17.de-de.xlf"
$ws.Range("H2").Value = "2016-08-30 13:04:37"

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$I$2') {
        $hl.Delete()
    }
}
$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"

$ws.Range("J2").Value = ""
$ws.Range("J2").Style = "Normal"

$ws.Range("K2").Value = "0001-01-01 00:00:00"

$ws.Columns.Item(9).ColumnWidth = 17.83
$ws.Columns.Item(10).ColumnWidth = 20.85
